# Auto-committed on 2022/04/07 週四
# Adds a new "OpDate" (作業日期) field to the CustRmk table schema (DBD sheet)
# and updates the related key/order-condition metadata on the DBS sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "DBD": table schema listing
# ---------------------------------------------------------------------
$dbd = $wb.Worksheets.Item("DBD")

# PrimaryKey cell now also includes OpDate
$dbd.Range("C3").Value = "CustNo,OpDate,RmkNo"

# Insert a new schema row for the OpDate field right after the CustNo row
$dbd.Rows.Item(10).Insert()

# Carry the surrounding table formatting (borders, font, ...) down into the
# freshly inserted row before overwriting it with the OpDate field values
$dbd.Range("A9:G9").Copy()
$dbd.Range("A10:G10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$dbd.Range("A10").Value = 2
$dbd.Range("B10").Value = "OpDate"
$dbd.Range("C10").Value = "作業日期"
$dbd.Range("D10").Value = "DECIMALD"
$dbd.Range("E10").Value = 8
$dbd.Range("G10").Value = "2022.4.7 by eric
轉換原日期"

# Formatting to match the rest of the table plus a highlight fill + comment wrap
$newRowRange = $dbd.Range("A10:G10")
$newRowRange.Interior.Color = 65535
$newRowRange.HorizontalAlignment = -4131
$newRowRange.VerticalAlignment = -4160
$newRowRange.WrapText = $true
$dbd.Rows.Item(10).RowHeight = 32.4

# ---------------------------------------------------------------------
# Sheet "DBS": key / order condition lookup table
# ---------------------------------------------------------------------
$dbs = $wb.Worksheets.Item("DBS")

# maxRmkNoFirst now also filters on OpDate
$dbs.Range("B4").Value = "CustNo = ,AND OpDate ="

# New "其他ORDER條件" for findCustNo
$dbs.Range("C2").Value = "OpDate,RmkNo"

Write-Output "edit complete"
